$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.558.91"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +0.35%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.560.23"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +0.60%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "606.86"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.28%  "

$ws.Range("E6").Value = "  +0.59%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.558.72"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.60%  "

$ws.Range("E8").Value = "  +0.06%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.497"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +3.45%  "

$ws.Range("E10").Value = "  -0.60%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.99"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -2.02%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.414"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +0.76%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.163.79"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +0.60%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000207"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -0.32%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "30.03"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -0.84%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.559.65"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +0.62%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "66.638.38"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +0.34%  "

$ws.Range("E18").Value = "  +0.33%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.54"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +5.50%  "

$ws.Range("E20").Value = "  +0.08%  "

$ws.Range("E21").Value = "  +0.26%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "431.40"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +1.15%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "79.47"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.70%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.700.33"
$ws.Range("D25").ClearFormats()

$ws.Range("E26").Value = "  -0.18%  "

$ws.Range("E27").Value = "  -0.39%  "

$ws.Range("E28").Value = "  -1.82%  "

$ws.Range("E29").Value = "  +1.06%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.16"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -0.95%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.00"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.02%  "

$ws.Range("E32").Value = "  -2.28%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.555.26"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.73%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "25.32"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +0.08%  "

$ws.Range("E35").Value = "  -3.87%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "7.83"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +0.19%  "

$ws.Range("E37").Value = "  +0.02%  "

$ws.Range("E38").Value = "  -1.84%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.60"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.42%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "174.00"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +0.34%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0847"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -1.16%  "

$ws.Range("E42").Value = "  -1.42%  "

$ws.Range("E43").Value = "  -0.66%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.93"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +1.37%  "

$ws.Range("E45").Value = "  +0.07%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.51"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +4.40%  "

$ws.Range("E47").Value = "  -2.74%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "25.07"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -3.91%  "

$ws.Range("E49").Value = "  +0.55%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "23.50"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +3.74%  "

$ws.Range("E51").Value = "  -0.27%  "
